$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace author first/last names with corrected data
$ws.Range("A1").Value = "Olivia"
$ws.Range("B1").Value = "Benet"
$ws.Range("A2").Value = "Ethan"
$ws.Range("B2").Value = "Harris"
$ws.Range("A3").Value = "Lily"
$ws.Range("B3").Value = "Thomson"
$ws.Range("A4").Value = "James"
$ws.Range("B4").Value = "Walker"
$ws.Range("A5").Value = "Liam"
$ws.Range("B5").Value = "Edwards"

# Move the active selection to D4
$ws.Range("D4").Select()
